$wb = $excel.ActiveWorkbook

$wsHCAHeart  = $wb.Worksheets.Item("HCAHeart_ID")
$wsControls  = $wb.Worksheets.Item("Controls_Synonymous")

# The edit happened on the "Controls_Synonymous" sheet: cell B6 ("H7") was
# corrected to "H6", and the in-sheet selection moved on to B7.
$wsControls.Activate()
$wsControls.Range("B6").Value = "H6"
$wsControls.Range("B7").Select()

# Afterwards the workbook was left with the first sheet ("HCAHeart_ID")
# as the active/selected tab (its own selection, C24, is unchanged).
$wsHCAHeart.Activate()

# Reposition the application window to match the saved view state.
$win = $excel.ActiveWindow
$win.Left = 30320
$win.Top = 1180
